$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting
# Late / heading / Outstanding one column to the right (N -> O, O -> P, P -> Q).
$ws.Columns("N").Insert()

# Match the width Excel gives the freshly inserted column (it inherits the
# width of the column immediately to its left, "In Advance" = 10.7109375).
$ws.Columns("N").ColumnWidth = 9.86

# Make "Repayment schedule" the active sheet/tab with J17 selected.
$ws.Activate() | Out-Null
$ws.Range("J17").Select() | Out-Null
